$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 51
$ws.Range("H51").Value = 8450.5
$ws.Range("I51").Value = 13950
$ws.Range("J51").Value = 2951
$ws.Range("K51").Value = 13950
$ws.Range("L51").Value = 2951
$ws.Range("M51").Value = -13466
$ws.Range("N51").Value = -3919

# Row 76
$ws.Range("H76").Value = 4424.75
$ws.Range("I76").Value = 4403
$ws.Range("J76").Value = 4432
$ws.Range("K76").Value = 4403
$ws.Range("L76").Value = 4432
$ws.Range("M76").Value = -4088
$ws.Range("N76").Value = -5062

# Row 79
$ws.Range("H79").Value = 4424.75
$ws.Range("I79").Value = 4403
$ws.Range("J79").Value = 4432
$ws.Range("K79").Value = 4403
$ws.Range("L79").Value = 4432
$ws.Range("M79").Value = -3311
$ws.Range("N79").Value = -6616

# Row 80
$ws.Range("H80").Value = 250045630
$ws.Range("I80").Value = 500001250
$ws.Range("J80").Value = 90003
$ws.Range("K80").Value = 1500003750
$ws.Range("L80").Value = 270009
$ws.Range("M80").Value = -1500002752
$ws.Range("N80").Value = -272005

# Row 83
$ws.Range("H83").Value = 250045630
$ws.Range("I83").Value = 500001250
$ws.Range("J83").Value = 90003
$ws.Range("K83").Value = 4500011250
$ws.Range("L83").Value = 810027
$ws.Range("M83").Value = -4500006258
$ws.Range("N83").Value = -820011

# Row 100
$ws.Range("H100").Value = 3459.15
$ws.Range("I100").Value = 2745.3333
$ws.Range("J100").Value = 3765.0715
$ws.Range("K100").Value = 2745.3333
$ws.Range("L100").Value = 3765.0715
$ws.Range("M100").Value = -2204.3333
$ws.Range("N100").Value = -4847.0715

# Row 106
$ws.Range("H106").Value = 6063147
$ws.Range("I106").Value = 6669361.5
$ws.Range("J106").Value = 1000
$ws.Range("K106").Value = 6669361.5
$ws.Range("L106").Value = 1000
$ws.Range("M106").Value = -6668730.5

# Row 137
$ws.Range("H137").Value = 39149.363
$ws.Range("I137").Value = 42889.75
$ws.Range("J137").Value = 1745.5
$ws.Range("K137").Value = 128669.25
$ws.Range("L137").Value = 5236.5
$ws.Range("M137").Value = -126119.25
$ws.Range("N137").Value = -10336.5

# Row 138
$ws.Range("H138").Value = 3805.36
$ws.Range("I138").Value = 2739.9443
$ws.Range("J138").Value = 4039.2317
$ws.Range("K138").Value = 8219.832900000001
$ws.Range("L138").Value = 12117.6951
$ws.Range("M138").Value = -3079.832900000001
$ws.Range("N138").Value = -22397.6951

# Row 141
$ws.Range("H141").Value = 919.8
$ws.Range("I141").Value = 919.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2759.4
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2420.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 15875651
$ws.Range("I32").Value = 18869972
$ws.Range("J32").Value = 5750.6
$ws.Range("K32").Value = 18869972
$ws.Range("L32").Value = 5750.6
$ws.Range("M32").Value = -18869685
$ws.Range("N32").Value = -6324.6

# Row 45
$ws.Range("H45").Value = 2823.0476
$ws.Range("I45").Value = 2569.4546
$ws.Range("J45").Value = 3752.889
$ws.Range("K45").Value = 2569.4546
$ws.Range("L45").Value = 3752.889
$ws.Range("M45").Value = -2192.4546
$ws.Range("N45").Value = -4506.889

# Row 61
$ws.Range("H61").Value = 2873.4736
$ws.Range("I61").Value = 2828.2856
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2828.2856
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2616.2856
$ws.Range("N61").Value = -3424

# Row 110
$ws.Range("M110").ClearContents()
$ws.Range("H110").Value = 3500
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 3500
$ws.Range("N110").Value = -7590

# Row 117
$ws.Range("H117").Value = 86527.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 86527.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 86527.5
$ws.Range("N117").Value = -95705.5

# Row 118
$ws.Range("H118").Value = 118665.664
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 118665.664
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 118665.664
$ws.Range("N118").Value = -121979.664

# Row 136
$ws.Range("H136").Value = 2873.4736
$ws.Range("I136").Value = 2828.2856
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 8484.856800000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -5934.856800000001
$ws.Range("N136").Value = -14100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 1906.4584
$ws.Range("I99").Value = 1339.8125
$ws.Range("J99").Value = 3039.75
$ws.Range("K99").Value = 1339.8125
$ws.Range("L99").Value = 3039.75
$ws.Range("M99").Value = 158.1875
$ws.Range("N99").Value = -6035.75

# Row 116
$ws.Range("H116").Value = 96244
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 96244
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 96244
$ws.Range("N116").Value = -105422

# Row 141
$ws.Range("H141").Value = 188249.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 188249.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 188249.5
$ws.Range("N141").Value = -198609.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 9290.565000000001
$ws.Range("I31").Value = 2384.375
$ws.Range("J31").Value = 12973.866
$ws.Range("K31").Value = 2384.375
$ws.Range("L31").Value = 12973.866
$ws.Range("M31").Value = -2089.375
$ws.Range("N31").Value = -13563.866

# Row 34
$ws.Range("H34").Value = 9290.565000000001
$ws.Range("I34").Value = 2384.375
$ws.Range("J34").Value = 12973.866
$ws.Range("K34").Value = 2384.375
$ws.Range("L34").Value = 12973.866
$ws.Range("M34").Value = -2182.375
$ws.Range("N34").Value = -13377.866

# Row 105
$ws.Range("H105").Value = 1784.7727
$ws.Range("I105").Value = 1071
$ws.Range("J105").Value = 3314.2856
$ws.Range("K105").Value = 1071
$ws.Range("L105").Value = 3314.2856
$ws.Range("M105").Value = 676

# Row 117
$ws.Range("H117").Value = 51699.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 51699.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 51699.5
$ws.Range("N117").Value = -60877.5

# Row 118
$ws.Range("H118").Value = 135999.5
$ws.Range("I118").Value = 119000
$ws.Range("J118").Value = 152999
$ws.Range("K118").Value = 119000
$ws.Range("L118").Value = 152999
$ws.Range("M118").Value = -117343
$ws.Range("N118").Value = -156313

# Row 122
$ws.Range("H122").Value = 5006081.5
$ws.Range("I122").Value = 6672112
$ws.Range("J122").Value = 7990
$ws.Range("K122").Value = 20016336
$ws.Range("L122").Value = 23970
$ws.Range("M122").Value = -20013886

# Row 132
$ws.Range("H132").Value = 3477.8076
$ws.Range("I132").Value = 4384.4116
$ws.Range("J132").Value = 1765.3334
$ws.Range("K132").Value = 13153.2348
$ws.Range("L132").Value = 5296.0002
$ws.Range("M132").Value = -10623.2348

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 138
$ws.Range("H138").Value = 23668222
$ws.Range("I138").Value = 1999.2858
$ws.Range("J138").Value = 106500000
$ws.Range("K138").Value = 5997.857400000001
$ws.Range("L138").Value = 319500000
$ws.Range("M138").Value = -857.8574000000008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 11
$ws.Range("H11").Value = 5203800
$ws.Range("I11").Value = 20000000
$ws.Range("J11").Value = 1504750
$ws.Range("K11").Value = 20000000
$ws.Range("L11").Value = 1504750
$ws.Range("M11").Value = -19999861
$ws.Range("N11").Value = -1505028

# Row 132
$ws.Range("H132").Value = 4224.478
$ws.Range("I132").Value = 3558.5293
$ws.Range("J132").Value = 6111.3335
$ws.Range("K132").Value = 10675.5879
$ws.Range("L132").Value = 18334.0005
$ws.Range("M132").Value = -8145.5879

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 3852.75
$ws.Range("I22").Value = 3483
$ws.Range("J22").Value = 4074.6
$ws.Range("K22").Value = 3483
$ws.Range("L22").Value = 4074.6
$ws.Range("M22").Value = -3188
$ws.Range("N22").Value = -4664.6

# Row 27
$ws.Range("H27").Value = 3852.75
$ws.Range("I27").Value = 3483
$ws.Range("J27").Value = 4074.6
$ws.Range("K27").Value = 3483
$ws.Range("L27").Value = 4074.6
$ws.Range("M27").Value = -3376
$ws.Range("N27").Value = -4288.6

# Row 122
$ws.Range("H122").Value = 7421.4546
$ws.Range("I122").Value = 5070.778
$ws.Range("J122").Value = 17999.5
$ws.Range("K122").Value = 15212.334
$ws.Range("L122").Value = 53998.5
$ws.Range("M122").Value = -12762.334

# Row 136
$ws.Range("H136").Value = 1496.32
$ws.Range("I136").Value = 924.2381
$ws.Range("J136").Value = 4499.75
$ws.Range("K136").Value = 2772.7143
$ws.Range("L136").Value = 13499.25
$ws.Range("M136").Value = -222.7143000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 39
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 15044
$ws.Range("I39").Value = 15044
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 15044
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -14631

# Row 42
$ws.Range("H42").Value = 21461
$ws.Range("I42").Value = 21461
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 21461
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -21083

# Row 43
$ws.Range("H43").Value = 25261
$ws.Range("I43").Value = 25261
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 25261
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -25112

# Row 81
$ws.Range("H81").Value = 3732.6155
$ws.Range("I81").Value = 2816.25
$ws.Range("J81").Value = 5198.8
$ws.Range("K81").Value = 5632.5
$ws.Range("L81").Value = 10397.6
$ws.Range("M81").Value = -4571.5

# Row 84
$ws.Range("H84").Value = 3732.6155
$ws.Range("I84").Value = 2816.25
$ws.Range("J84").Value = 5198.8
$ws.Range("K84").Value = 28162.5
$ws.Range("L84").Value = 51988
$ws.Range("M84").Value = -22858.5

# Row 113
$ws.Range("H113").Value = 752.1429000000001
$ws.Range("I113").Value = 894.1
$ws.Range("J113").Value = 397.25
$ws.Range("K113").Value = 2682.3
$ws.Range("L113").Value = 1191.75
$ws.Range("M113").Value = -512.3000000000002
$ws.Range("N113").Value = -5531.75

# Row 125
$ws.Range("H125").Value = 56142.855
$ws.Range("I125").Value = 43000
$ws.Range("J125").Value = 58333.332
$ws.Range("K125").Value = 43000
$ws.Range("L125").Value = 58333.332
$ws.Range("M125").Value = -38080
$ws.Range("N125").Value = -68173.33199999999

# Row 126
$ws.Range("H126").Value = 3628.6667
$ws.Range("I126").Value = 3759.3572
$ws.Range("J126").Value = 1799
$ws.Range("K126").Value = 11278.0716
$ws.Range("L126").Value = 5397
$ws.Range("M126").Value = -8808.071599999999
$ws.Range("N126").Value = -10337

# Row 136
$ws.Range("H136").Value = 30023.834
$ws.Range("I136").Value = 1633.6666
$ws.Range("J136").Value = 115194.336
$ws.Range("K136").Value = 4900.9998
$ws.Range("L136").Value = 345583.008
$ws.Range("M136").Value = -2350.9998

Write-Output "Applied all changes"